$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Labor Budget" row for the JONAS G54 source cell that now
# feeds the consolidation.
$ws.Range("A22").Value = "JONAS"
$ws.Range("B22").Value = "G54"
$ws.Range("C22").Value = "Labor Budget"

# Give the new row the same look as the rest of the data rows (thin
# border, centered, regular font/fill).
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null

# Row 2 loses its old bold/filled look and now matches the plain data rows.
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A2:C2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Leave the selection where it would land after entering the new row.
$ws.Range("A23").Select() | Out-Null
